$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-10-30", 7.4332848),
    @("2023-10-31", 23.751),
    @("2023-11-01", 136.578),
    @("2023-11-02", 49.40425),
    @("2023-11-03", 57.2145)
)

$startRow = 210
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $cellA = $ws.Cells.Item($row, 1)
    # Force the date-like string to be stored as literal text (matching
    # the source row format) instead of being auto-parsed into a date
    # serial number, then drop back to the default "Normal" style so no
    # stray cell-level number format is left behind.
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[$i][0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
